$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values.
# A leading backtick-apostrophe forces Excel to store the value as literal text
# (mirrors the source workbook, where these are inline/shared text strings, not
# numbers) even for values that otherwise look numeric (e.g. "22.27").
# Re-applying the "Normal" style afterwards clears the quote-prefix cell format
# so no stray style index is left attached to the cell.

$ws.Range("D2").Value = "`'27.654.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "`'  +0.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "`'1.589.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "`'  -0.23%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "`'  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "`'207.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "`'  +0.28%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "`'  +0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "`'  +0.11%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "`'22.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "`'  +0.19%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "`'  -0.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "`'  +0.20%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "`'  -0.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "`'1.815.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "`'  -0.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "`'1.568.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "`'  -1.51%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "`'  -0.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "`'0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "`'  -1.83%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "`'27.662.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "`'  +0.47%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "`'  -0.29%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "`'216.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "`'  +0.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "`'  +0.32%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "`'  -0.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "`'  +0.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "`'  -1.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "`'  +1.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "`'  -0.53%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "`'153.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "`'  -1.13%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "`'6.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "`'  +4.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "`'  +0.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "`'  +0.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "`'  -0.50%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "`'  -0.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "`'  +1.30%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "`'  -2.12%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "`'1.371.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "`'  +1.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "`'2.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "`'  +0.74%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "`'  -0.15%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "`'0.965"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "`'  +0.86%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "`'  -0.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "`'0.0169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "`'  +1.97%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "`'0.534"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "`'  -0.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "`'  +1.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "`'  +0.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "`'0.974"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "`'  +1.64%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "`'64.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "`'  +0.70%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "`'  +4.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "`'  +1.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "`'  -1.60%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "`'1.726.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "`'85.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "`'  -1.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "`'  +0.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "`'  -0.33%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "`'0.0495"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "`'  -0.45%  "
$ws.Range("E51").Style = "Normal"
